# Update countries & provincias Spain
# Applies the COVID-19 stats refresh: new case counts for a set of rows
# (which re-sorts a few country pairs by "Casos totales") and bumps the
# "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refreshed "Datos actualizados..." timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 02:33"

# --- Country names that swapped rank (column A) because their updated ---
# --- "Casos totales" (column B) changed their sort position ---
$ws.Range("A18").Value = "Argentina"
$ws.Range("A19").Value = "Banglades"

$ws.Range("A154").Value = "Bahamas"
$ws.Range("A155").Value = "Jamaica"

$ws.Range("A179").Value = "Eritrea"
$ws.Range("A180").Value = "Camboya"
$ws.Range("A181").Value = "Papua Nueva Guinea"

$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, ---
# --- Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5414450
$ws.Range("C4").Value = 54148
$ws.Range("D4").Value = 2833290
$ws.Range("E4").Value = 2410800
$ws.Range("G4").Value = 1229
$ws.Range("H4").Value = 170360

# Row 5: Brasil
$ws.Range("B5").Value = 3229621
$ws.Range("C5").Value = 59147
$ws.Range("E5").Value = 767417
$ws.Range("G5").Value = 1301
$ws.Range("H5").Value = 105564

# Row 18: Argentina (after swap)
$ws.Range("B18").Value = 276072
$ws.Range("C18").Value = 7498
$ws.Range("D18").Value = 192434
$ws.Range("E18").Value = 78276
$ws.Range("G18").Value = 149
$ws.Range("H18").Value = 5362

# Row 19: Banglades (after swap)
$ws.Range("B19").Value = 269115
$ws.Range("C19").Value = 2617
$ws.Range("D19").Value = 154871
$ws.Range("E19").Value = 110687
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 3557

# Row 27
$ws.Range("B27").Value = 121234
$ws.Range("C27").Value = 390
$ws.Range("D27").Value = 107553
$ws.Range("E27").Value = 4666
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 9015

# Row 39
$ws.Range("B39").Value = 78446
$ws.Range("C39").Value = 1069
$ws.Range("D39").Value = 52210
$ws.Range("E39").Value = 24514
$ws.Range("G39").Value = 19
$ws.Range("H39").Value = 1722

# Row 142
$ws.Range("B142").Value = 1409
$ws.Range("C142").Value = 16
$ws.Range("D142").Value = 1180
$ws.Range("E142").Value = 192

# Row 154: Bahamas (after swap)
$ws.Range("B154").Value = 1089
$ws.Range("C154").Value = 53
$ws.Range("D154").Value = 138
$ws.Range("E154").Value = 936
$ws.Range("H154").Value = 15

# Row 155: Jamaica (after swap)
$ws.Range("B155").Value = 1065
$ws.Range("C155").Value = 18
$ws.Range("D155").Value = 753
$ws.Range("E155").Value = 298
$ws.Range("H155").Value = 14

# Row 161
$ws.Range("B161").Value = 883
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 808

# Row 179: Eritrea (after rotation)
$ws.Range("B179").Value = 285
$ws.Range("D179").Value = 248
$ws.Range("E179").Value = 37
$ws.Range("H179").Value = 0

# Row 180: Camboya (after rotation)
$ws.Range("B180").Value = 272
$ws.Range("C180").Value = 4
$ws.Range("D180").Value = 223
$ws.Range("E180").Value = 49

# Row 181: Papua Nueva Guinea (after rotation)
$ws.Range("B181").Value = 271
$ws.Range("C181").Value = 2
$ws.Range("D181").Value = 78
$ws.Range("E181").Value = 190
$ws.Range("H181").Value = 3

# Row 188
$ws.Range("D188").Value = 118
$ws.Range("E188").Value = 19

# Row 213: Montserrat (after swap)
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214: Islas Malvinas (after swap)
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
